$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the "before" values of the columns that change (D, L, M, N, O, P, Q, S, T)
# for rows 2..11, then re-distribute them to the rows according to the new order
# observed in the updated dataset (a re-shuffle of the weekly price records).
$rows = 2..11
$before = @{}
foreach ($r in $rows) {
    $before[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        N = $ws.Cells.Item($r, 14).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
        Q = $ws.Cells.Item($r, 17).Value2
        S = $ws.Cells.Item($r, 19).Value2
        T = $ws.Cells.Item($r, 20).Value2
    }
}

# Mapping: new row -> source (old) row whose data now lands there.
$map = @{
    2  = 9
    3  = 5
    4  = 7
    5  = 10
    6  = 11
    7  = 8
    8  = 4
    9  = 6
    10 = 2
    11 = 3
}

foreach ($r in $rows) {
    $src = $before[$map[$r]]
    $ws.Cells.Item($r, 4).Value2 = $src.D
    $ws.Cells.Item($r, 12).Value2 = $src.L
    $ws.Cells.Item($r, 13).Value2 = $src.M
    $ws.Cells.Item($r, 14).Value2 = $src.N
    $ws.Cells.Item($r, 15).Value2 = $src.O
    $ws.Cells.Item($r, 16).Value2 = $src.P
    $ws.Cells.Item($r, 17).Value2 = $src.Q
    $ws.Cells.Item($r, 19).Value2 = $src.S
    $ws.Cells.Item($r, 20).Value2 = $src.T
}
